$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 2911.8  # H4: was 1139.5
$ws.Cells.Item(4, 9).Value = 3786.6667  # I4: was 679.5
$ws.Cells.Item(4, 11).Value = 3786.6667  # K4: was 679.5
$ws.Cells.Item(4, 13).Value = -3672.6667  # M4: was -565.5

$ws.Cells.Item(29, 8).Value = 1644.7273  # H29: was 2083.6667
$ws.Cells.Item(29, 9).Value = 219.8  # I29: was 249.5
$ws.Cells.Item(29, 10).Value = 2832.1667  # J29: was 3000.75
$ws.Cells.Item(29, 11).Value = 659.4000000000001  # K29: was 748.5
$ws.Cells.Item(29, 12).Value = 8496.500100000001  # L29: was 9002.25
$ws.Cells.Item(29, 13).Value = -378.4000000000001  # M29: was -467.5
$ws.Cells.Item(29, 14).Value = -9058.500100000001  # N29: was -9564.25

$ws.Cells.Item(70, 8).Value = 14999.667  # H70: was 15000
$ws.Cells.Item(70, 9).Value = 14999.667  # I70: was 15000
$ws.Cells.Item(70, 11).Value = 44999.001  # K70: was 45000
$ws.Cells.Item(70, 13).Value = -44729.001  # M70: was -44730

$ws.Cells.Item(73, 8).Value = 14999.667  # H73: was 15000
$ws.Cells.Item(73, 9).Value = 14999.667  # I73: was 15000
$ws.Cells.Item(73, 11).Value = 44999.001  # K73: was 45000
$ws.Cells.Item(73, 13).Value = -44063.001  # M73: was -44064

$ws.Cells.Item(116, 8).Value = 3631.3103  # H116: was 3753.9614
$ws.Cells.Item(116, 9).Value = 2947.5334  # I116: was 3007.923
$ws.Cells.Item(116, 10).Value = 4363.9287  # J116: was 4500
$ws.Cells.Item(116, 11).Value = 2947.5334  # K116: was 3007.923
$ws.Cells.Item(116, 12).Value = 4363.9287  # L116: was 4500
$ws.Cells.Item(116, 13).Value = 494.4666000000002  # M116: was 434.0770000000002
$ws.Cells.Item(116, 14).Value = -11247.9287  # N116: was -11384

$ws.Cells.Item(125, 8).Value = 2947  # H125: was 2878
$ws.Cells.Item(125, 10).Value = 2327  # J125: was 2326.5
$ws.Cells.Item(125, 12).Value = 20943  # L125: was 20938.5
$ws.Cells.Item(125, 14).Value = -25863  # N125: was -25858.5

$ws.Cells.Item(132, 8).Value = 2377.077  # H132: was 2638.1667
$ws.Cells.Item(132, 9).Value = 2325.1667  # I132: was 2505.8
$ws.Cells.Item(132, 10).Value = 3000  # J132: was 3300
$ws.Cells.Item(132, 11).Value = 6975.500100000001  # K132: was 7517.400000000001
$ws.Cells.Item(132, 12).Value = 9000  # L132: was 9900
$ws.Cells.Item(132, 13).Value = -4445.500100000001  # M132: was -4987.400000000001
$ws.Cells.Item(132, 14).Value = -14060  # N132: was -14960

$ws.Cells.Item(137, 8).Value = 905.4286  # H137: was 1245.5714
$ws.Cells.Item(137, 10).Value = 861.3333  # J137: was 1275.8
$ws.Cells.Item(137, 12).Value = 2583.9999  # L137: was 3827.4
$ws.Cells.Item(137, 14).Value = -7683.9999  # N137: was -8927.4

$ws.Cells.Item(141, 8).Value = 3055.1667  # H141: was 3267
$ws.Cells.Item(141, 9).Value = 2122.75  # I141: was 2165
$ws.Cells.Item(141, 11).Value = 6368.25  # K141: was 6495
$ws.Cells.Item(141, 13).Value = -1188.25  # M141: was -1315

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 572.2  # H5: was 622.2
$ws.Cells.Item(5, 9).Value = 465.25  # I5: was 527.75
$ws.Cells.Item(5, 11).Value = 465.25  # K5: was 527.75
$ws.Cells.Item(5, 13).Value = -353.25  # M5: was -415.75

$ws.Cells.Item(12, 8).Value = 7500  # H12: was 3
$ws.Cells.Item(12, 9).Value = 0  # I12: was 3
$ws.Cells.Item(12, 10).Value = 7500  # J12: was 0
$ws.Cells.Item(12, 11).Value = 0  # K12: was 3
$ws.Cells.Item(12, 12).ClearContents()  # L12: was 0
$ws.Cells.Item(12, 13).Value = 7500  # M12: was 170
$ws.Cells.Item(12, 14).Value = -7846  # N12: was None

$ws.Cells.Item(32, 8).Value = 2084113.2  # H32: was 2142012.8
$ws.Cells.Item(32, 9).Value = 1947560.9  # I32: was 2003213
$ws.Cells.Item(32, 11).Value = 1947560.9  # K32: was 2003213
$ws.Cells.Item(32, 13).Value = -1947273.9  # M32: was -2002926

$ws.Cells.Item(45, 8).Value = 5110.9  # H45: was 5634.4443
$ws.Cells.Item(45, 10).Value = 399  # J45: was 0
$ws.Cells.Item(45, 12).Value = 399  # L45: was 0
$ws.Cells.Item(45, 14).Value = -1153  # N45: was None

$ws.Cells.Item(61, 8).Value = 5544.9585  # H61: was 8773.913
$ws.Cells.Item(61, 9).Value = 5635.409  # I61: was 8809.091
$ws.Cells.Item(61, 10).Value = 4550  # J61: was 8000
$ws.Cells.Item(61, 11).Value = 5635.409  # K61: was 8809.091
$ws.Cells.Item(61, 12).Value = 4550  # L61: was 8000
$ws.Cells.Item(61, 13).Value = -5423.409  # M61: was -8597.091
$ws.Cells.Item(61, 14).Value = -4974  # N61: was -8424

$ws.Cells.Item(74, 8).Value = 1009.8  # H74: was 988.6
$ws.Cells.Item(74, 9).Value = 987.25  # I74: was 988.6
$ws.Cells.Item(74, 10).Value = 1100  # J74: was 0
$ws.Cells.Item(74, 11).Value = 987.25  # K74: was 988.6
$ws.Cells.Item(74, 12).Value = 1100  # L74: was 0
$ws.Cells.Item(74, 13).Value = -113.25  # M74: was -114.6
$ws.Cells.Item(74, 14).Value = -2848  # N74: was None

$ws.Cells.Item(77, 8).Value = 1009.8  # H77: was 988.6
$ws.Cells.Item(77, 9).Value = 987.25  # I77: was 988.6
$ws.Cells.Item(77, 10).Value = 1100  # J77: was 0
$ws.Cells.Item(77, 11).Value = 4936.25  # K77: was 4943
$ws.Cells.Item(77, 12).Value = 5500  # L77: was 0
$ws.Cells.Item(77, 13).Value = -568.25  # M77: was -575
$ws.Cells.Item(77, 14).Value = -14236  # N77: was None

$ws.Cells.Item(97, 8).Value = 793.9091  # H97: was 1047.7142
$ws.Cells.Item(97, 9).Value = 836.625  # I97: was 1132.8
$ws.Cells.Item(97, 10).Value = 680  # J97: was 835
$ws.Cells.Item(97, 11).Value = 836.625  # K97: was 1132.8
$ws.Cells.Item(97, 12).Value = 680  # L97: was 835
$ws.Cells.Item(97, 13).Value = -340.625  # M97: was -636.8
$ws.Cells.Item(97, 14).Value = -1672  # N97: was -1827

$ws.Cells.Item(122, 8).Value = 17198.95  # H122: was 18041.21
$ws.Cells.Item(122, 9).Value = 17198.95  # I122: was 18041.21
$ws.Cells.Item(122, 11).Value = 51596.85000000001  # K122: was 54123.63
$ws.Cells.Item(122, 13).Value = -49146.85000000001  # M122: was -51673.63

$ws.Cells.Item(132, 8).Value = 2194.5715  # H132: was 2194.8
$ws.Cells.Item(132, 9).Value = 2191.5  # I132: was 2192
$ws.Cells.Item(132, 10).Value = 2198.6667  # J132: was 2199
$ws.Cells.Item(132, 11).Value = 6574.5  # K132: was 6576
$ws.Cells.Item(132, 12).Value = 6596.000100000001  # L132: was 6597
$ws.Cells.Item(132, 13).Value = -4044.5  # M132: was -4046
$ws.Cells.Item(132, 14).Value = -11656.0001  # N132: was -11657

$ws.Cells.Item(136, 8).Value = 5544.9585  # H136: was 8773.913
$ws.Cells.Item(136, 9).Value = 5635.409  # I136: was 8809.091
$ws.Cells.Item(136, 10).Value = 4550  # J136: was 8000
$ws.Cells.Item(136, 11).Value = 16906.227  # K136: was 26427.273
$ws.Cells.Item(136, 12).Value = 13650  # L136: was 24000
$ws.Cells.Item(136, 13).Value = -14356.227  # M136: was -23877.273
$ws.Cells.Item(136, 14).Value = -18750  # N136: was -29100

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 572.2  # H4: was 622.2
$ws.Cells.Item(4, 9).Value = 465.25  # I4: was 527.75
$ws.Cells.Item(4, 11).Value = 465.25  # K4: was 527.75
$ws.Cells.Item(4, 13).Value = -350.25  # M4: was -412.75

$ws.Cells.Item(134, 8).Value = 4984.5  # H134: was 6147.25
$ws.Cells.Item(134, 9).Value = 5016.8  # I134: was 6297
$ws.Cells.Item(134, 11).Value = 15050.4  # K134: was 18891
$ws.Cells.Item(134, 13).Value = -12515.4  # M134: was -16356

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1998.6666  # H22: was 2099.3333
$ws.Cells.Item(22, 9).Value = 1998  # I22: was 2149
$ws.Cells.Item(22, 10).Value = 1999  # J22: was 2000
$ws.Cells.Item(22, 11).Value = 1998  # K22: was 2149
$ws.Cells.Item(22, 12).Value = 1999  # L22: was 2000
$ws.Cells.Item(22, 13).Value = -1648  # M22: was -1799
$ws.Cells.Item(22, 14).Value = -2699  # N22: was -2700

$ws.Cells.Item(31, 8).Value = 1431.7894  # H31: was 1408.1578
$ws.Cells.Item(31, 9).Value = 1727  # I31: was 1396.1666
$ws.Cells.Item(31, 10).Value = 1353.0667  # J31: was 1413.6923
$ws.Cells.Item(31, 11).Value = 1727  # K31: was 1396.1666
$ws.Cells.Item(31, 12).Value = 1353.0667  # L31: was 1413.6923
$ws.Cells.Item(31, 13).Value = -1432  # M31: was -1101.1666
$ws.Cells.Item(31, 14).Value = -1943.0667  # N31: was -2003.6923

$ws.Cells.Item(34, 8).Value = 1431.7894  # H34: was 1408.1578
$ws.Cells.Item(34, 9).Value = 1727  # I34: was 1396.1666
$ws.Cells.Item(34, 10).Value = 1353.0667  # J34: was 1413.6923
$ws.Cells.Item(34, 11).Value = 1727  # K34: was 1396.1666
$ws.Cells.Item(34, 12).Value = 1353.0667  # L34: was 1413.6923
$ws.Cells.Item(34, 13).Value = -1525  # M34: was -1194.1666
$ws.Cells.Item(34, 14).Value = -1757.0667  # N34: was -1817.6923

$ws.Cells.Item(35, 8).Value = 831.25  # H35: was 775
$ws.Cells.Item(35, 9).Value = 541.6667  # I35: was 312.5
$ws.Cells.Item(35, 11).Value = 541.6667  # K35: was 312.5
$ws.Cells.Item(35, 13).Value = -247.6667  # M35: was -18.5

$ws.Cells.Item(86, 8).Value = 6465.8887  # H86: was 6680.375
$ws.Cells.Item(86, 10).Value = 4816  # J86: was 4849
$ws.Cells.Item(86, 12).Value = 4816  # L86: was 4849
$ws.Cells.Item(86, 14).Value = -7062  # N86: was -7095

$ws.Cells.Item(89, 8).Value = 6465.8887  # H89: was 6680.375
$ws.Cells.Item(89, 10).Value = 4816  # J89: was 4849
$ws.Cells.Item(89, 12).Value = 24080  # L89: was 24245
$ws.Cells.Item(89, 14).Value = -35312  # N89: was -35477

$ws.Cells.Item(105, 8).Value = 3624.75  # H105: was 4250
$ws.Cells.Item(105, 9).Value = 1599.5  # I105: was 3500
$ws.Cells.Item(105, 10).Value = 4299.8335  # J105: was 4400
$ws.Cells.Item(105, 11).Value = 1599.5  # K105: was 3500
$ws.Cells.Item(105, 12).Value = 4299.8335  # L105: was 4400
$ws.Cells.Item(105, 13).Value = 147.5  # M105: was -1753
$ws.Cells.Item(105, 14).Value = -7793.8335  # N105: was -7894

$ws.Cells.Item(134, 8).Value = 2796.5  # H134: was 2799
$ws.Cells.Item(134, 9).Value = 2796.5  # I134: was 2799
$ws.Cells.Item(134, 11).Value = 8389.5  # K134: was 8397
$ws.Cells.Item(134, 13).Value = -5854.5  # M134: was -5862

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 9966259  # H4: was 8097608
$ws.Cells.Item(4, 9).Value = 12002956  # I4: was 9233070
$ws.Cells.Item(4, 11).Value = 36008868  # K4: was 27699210
$ws.Cells.Item(4, 13).Value = -36008756  # M4: was -27699098

$ws.Cells.Item(121, 8).Value = 10905.733  # H121: was 10535
$ws.Cells.Item(121, 9).Value = 27951  # I121: was 37086.332
$ws.Cells.Item(121, 10).Value = 4707.4546  # J121: was 4407.769
$ws.Cells.Item(121, 11).Value = 83853  # K121: was 111258.996
$ws.Cells.Item(121, 12).Value = 14122.3638  # L121: was 13223.307
$ws.Cells.Item(121, 13).Value = -82543  # M121: was -109948.996
$ws.Cells.Item(121, 14).Value = -16742.3638  # N121: was -15843.307

$ws.Cells.Item(131, 8).Value = 716795.0600000001  # H131: was 478646.44
$ws.Cells.Item(131, 10).Value = 716795.0600000001  # J131: was 478646.44
$ws.Cells.Item(131, 12).Value = 2150385.18  # L131: was 1435939.32
$ws.Cells.Item(131, 14).Value = -2160465.18  # N131: was -1446019.32

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 0  # H70: was 1494
$ws.Cells.Item(70, 9).Value = 0  # I70: was 1494
$ws.Cells.Item(70, 11).Value = 0  # K70: was 1494
$ws.Cells.Item(70, 13).ClearContents()  # M70: was -1224

$ws.Cells.Item(73, 8).Value = 0  # H73: was 1494
$ws.Cells.Item(73, 9).Value = 0  # I73: was 1494
$ws.Cells.Item(73, 11).Value = 0  # K73: was 1494
$ws.Cells.Item(73, 13).ClearContents()  # M73: was -558

$ws.Cells.Item(80, 8).Value = 4744.4  # H80: was 5586.75
$ws.Cells.Item(80, 9).Value = 3237.5  # I80: was 3733.3333
$ws.Cells.Item(80, 10).Value = 5749  # J80: was 6698.8
$ws.Cells.Item(80, 11).Value = 3237.5  # K80: was 3733.3333
$ws.Cells.Item(80, 12).Value = 5749  # L80: was 6698.8
$ws.Cells.Item(80, 13).Value = -2239.5  # M80: was -2735.3333
$ws.Cells.Item(80, 14).Value = -7745  # N80: was -8694.799999999999

$ws.Cells.Item(83, 8).Value = 4744.4  # H83: was 5586.75
$ws.Cells.Item(83, 9).Value = 3237.5  # I83: was 3733.3333
$ws.Cells.Item(83, 10).Value = 5749  # J83: was 6698.8
$ws.Cells.Item(83, 11).Value = 16187.5  # K83: was 18666.6665
$ws.Cells.Item(83, 12).Value = 28745  # L83: was 33494
$ws.Cells.Item(83, 13).Value = -11195.5  # M83: was -13674.6665
$ws.Cells.Item(83, 14).Value = -38729  # N83: was -43478

$ws.Cells.Item(113, 8).Value = 2679.5  # H113: was 2929.4285
$ws.Cells.Item(113, 10).Value = 5076.6665  # J113: was 7150
$ws.Cells.Item(113, 12).Value = 5076.6665  # L113: was 7150
$ws.Cells.Item(113, 14).Value = -9416.666499999999  # N113: was -11490

$ws.Cells.Item(122, 8).Value = 2799.2856  # H122: was 2999.25
$ws.Cells.Item(122, 9).Value = 2749.5  # I122: was 3466
$ws.Cells.Item(122, 10).Value = 2865.6667  # J122: was 2719.2
$ws.Cells.Item(122, 11).Value = 8248.5  # K122: was 10398
$ws.Cells.Item(122, 12).Value = 8597.000100000001  # L122: was 8157.599999999999
$ws.Cells.Item(122, 13).Value = -5798.5  # M122: was -7948
$ws.Cells.Item(122, 14).Value = -13497.0001  # N122: was -13057.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 669.8  # H16: was 801.3333
$ws.Cells.Item(16, 9).Value = 669.8  # I16: was 801.3333
$ws.Cells.Item(16, 11).Value = 669.8  # K16: was 801.3333
$ws.Cells.Item(16, 13).Value = -499.8  # M16: was -631.3333

$ws.Cells.Item(22, 8).Value = 3256.7144  # H22: was 3574.5
$ws.Cells.Item(22, 9).Value = 2550  # I22: was 2950
$ws.Cells.Item(22, 11).Value = 2550  # K22: was 2950
$ws.Cells.Item(22, 13).Value = -2255  # M22: was -2655

$ws.Cells.Item(27, 8).Value = 3256.7144  # H27: was 3574.5
$ws.Cells.Item(27, 9).Value = 2550  # I27: was 2950
$ws.Cells.Item(27, 11).Value = 2550  # K27: was 2950
$ws.Cells.Item(27, 13).Value = -2443  # M27: was -2843

$ws.Cells.Item(40, 8).Value = 2257.1667  # H40: was 2323.5557
$ws.Cells.Item(40, 9).Value = 1835.4  # I40: was 1925
$ws.Cells.Item(40, 10).Value = 2784.375  # J40: was 2821.75
$ws.Cells.Item(40, 11).Value = 1835.4  # K40: was 1925
$ws.Cells.Item(40, 12).Value = 2784.375  # L40: was 2821.75
$ws.Cells.Item(40, 13).Value = -1699.4  # M40: was -1789
$ws.Cells.Item(40, 14).Value = -3056.375  # N40: was -3093.75

$ws.Cells.Item(122, 8).Value = 5850.25  # H122: was 5780.39
$ws.Cells.Item(122, 9).Value = 4479.9443  # I122: was 4401.316
$ws.Cells.Item(122, 11).Value = 13439.8329  # K122: was 13203.948
$ws.Cells.Item(122, 13).Value = -10989.8329  # M122: was -10753.948

$ws.Cells.Item(132, 8).Value = 2193.7742  # H132: was 2251.5
$ws.Cells.Item(132, 9).Value = 2080.2  # I132: was 2147.625
$ws.Cells.Item(132, 11).Value = 6240.599999999999  # K132: was 6442.875
$ws.Cells.Item(132, 13).Value = -3710.599999999999  # M132: was -3912.875

$ws.Cells.Item(136, 8).Value = 2578.3635  # H136: was 2516.2
$ws.Cells.Item(136, 9).Value = 2499.6667  # I136: was 2350.75
$ws.Cells.Item(136, 10).Value = 2607.875  # J136: was 2626.5
$ws.Cells.Item(136, 11).Value = 7499.000100000001  # K136: was 7052.25
$ws.Cells.Item(136, 12).Value = 7823.625  # L136: was 7879.5
$ws.Cells.Item(136, 13).Value = -4949.000100000001  # M136: was -4502.25
$ws.Cells.Item(136, 14).Value = -12923.625  # N136: was -12979.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(119, 8).Value = 0  # H119: was 78924
$ws.Cells.Item(119, 10).Value = 0  # J119: was 78924
$ws.Cells.Item(119, 12).ClearContents()  # L119: was 78924
$ws.Cells.Item(119, 14).Value = 0  # N119: was -88600

$ws.Cells.Item(120, 8).Value = 0  # H120: was 117500
$ws.Cells.Item(120, 10).Value = 0  # J120: was 117500
$ws.Cells.Item(120, 12).ClearContents()  # L120: was 117500
$ws.Cells.Item(120, 14).Value = 0  # N120: was -127176

$ws.Cells.Item(122, 8).Value = 4364.8423  # H122: was 4366.05
$ws.Cells.Item(122, 9).Value = 3135.1538  # I122: was 3181.8572
$ws.Cells.Item(122, 10).Value = 7029.1665  # J122: was 7129.1665
$ws.Cells.Item(122, 11).Value = 9405.4614  # K122: was 9545.571599999999
$ws.Cells.Item(122, 12).Value = 21087.4995  # L122: was 21387.4995
$ws.Cells.Item(122, 13).Value = -6955.4614  # M122: was -7095.571599999999
$ws.Cells.Item(122, 14).Value = -25987.4995  # N122: was -26287.4995

$ws.Cells.Item(132, 8).Value = 2026.5  # H132: was 1900.2222
$ws.Cells.Item(132, 10).Value = 1499.5  # J132: was 1296.3334
$ws.Cells.Item(132, 12).Value = 4498.5  # L132: was 3889.0002
$ws.Cells.Item(132, 14).Value = -9558.5  # N132: was -8949.0002
